# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets to match the latest scraped data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1985
$ws1.Range("F4").Value = 844
$ws1.Range("F5").Value = 1021
$ws1.Range("F6").Value = 346

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1985
$ws4.Range("F5").Value = 844
$ws4.Range("F6").Value = 1021
$ws4.Range("F7").Value = 346
